$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 10, 11, 14, 15, 17, 23, 26, 30, 31, 38, 43, 45, 46, 47, 49, 53, 55, 59, 60, 64, 65, 68, 71, 73, 79, 85, 86, 89, 94, 98, 100, 101, 103, 108, 115, 116, 117, 120, 123, 125, 130, 132, 134, 135, 140, 146, 150, 153, 156, 158, 162, 163, 164, 167, 171, 177, 179, 181, 182, 185, 188, 191, 193)

foreach ($r in $rows) {
    $ws.Range("G" + $r + ":N" + $r).Value = "N/A"
}
